$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (G=5503)
$ws.Range("H5").Value = 71.166664
$ws.Range("I5").Value = 78.40000000000001
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 78.40000000000001
$ws.Range("L5").Value = 35
$ws.Range("M5").Value = 36.59999999999999
$ws.Range("N5").Value = -265

# Row 19 (G=7015)
$ws.Range("H19").Value = 12557.8125
$ws.Range("I19").Value = 90.77778000000001
$ws.Range("J19").Value = 28586.857
$ws.Range("K19").Value = 90.77778000000001
$ws.Range("L19").Value = 28586.857
$ws.Range("M19").Value = 84.22221999999999
$ws.Range("N19").Value = -28936.857

# Row 32 (G=5484)
$ws.Range("H32").Value = 2499.5
$ws.Range("I32").Value = 2499.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2499.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2173.5

# Row 43 (G=5472)
$ws.Range("H43").Value = 1766.6666
$ws.Range("I43").Value = 1720
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 1720
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -1651
$ws.Range("N43").Value = -2138

# Row 74 (G=5507)
$ws.Range("M74").ClearContents()
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 4500
$ws.Range("N74").Value = -6372

# Row 77 (G=5507)
$ws.Range("M77").ClearContents()
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 22500
$ws.Range("N77").Value = -31860

# Row 132 (G=44049)
$ws.Range("H132").Value = 4496.294
$ws.Range("I132").Value = 4250.522
$ws.Range("J132").Value = 6757.4
$ws.Range("K132").Value = 12751.566
$ws.Range("L132").Value = 20272.2
$ws.Range("M132").Value = -10221.566
$ws.Range("N132").Value = -25332.2

# Row 137 (G=44013)
$ws.Range("H137").Value = 6424.5557
$ws.Range("I137").Value = 7750
$ws.Range("J137").Value = 6045.857
$ws.Range("K137").Value = 23250
$ws.Range("L137").Value = 18137.571
$ws.Range("M137").Value = -20700
$ws.Range("N137").Value = -23237.571

# Row 138 (G=44169)
$ws.Range("H138").Value = 5515.4097
$ws.Range("I138").Value = 5019.2144
$ws.Range("J138").Value = 5663.213
$ws.Range("K138").Value = 15057.6432
$ws.Range("L138").Value = 16989.639
$ws.Range("M138").Value = -9917.643199999999
$ws.Range("N138").Value = -27269.639


$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G=44147)
$ws.Range("H32").Value = 20340.738
$ws.Range("I32").Value = 16776.982
$ws.Range("J32").Value = 45732.5
$ws.Range("K32").Value = 16776.982
$ws.Range("L32").Value = 45732.5
$ws.Range("M32").Value = -16489.982
$ws.Range("N32").Value = -46306.5

# Row 61 (G=43999)
$ws.Range("H61").Value = 1255560.1
$ws.Range("I61").Value = 3837.3333
$ws.Range("J61").Value = 2006593.8
$ws.Range("K61").Value = 3837.3333
$ws.Range("L61").Value = 2006593.8
$ws.Range("M61").Value = -3625.3333
$ws.Range("N61").Value = -2007017.8

# Row 74 (G=44000)
$ws.Range("H74").Value = 3187.6897
$ws.Range("I74").Value = 2628.0952
$ws.Range("J74").Value = 4656.625
$ws.Range("K74").Value = 2628.0952
$ws.Range("L74").Value = 4656.625
$ws.Range("M74").Value = -1754.0952
$ws.Range("N74").Value = -6404.625

# Row 77 (G=44000)
$ws.Range("H77").Value = 3187.6897
$ws.Range("I77").Value = 2628.0952
$ws.Range("J77").Value = 4656.625
$ws.Range("K77").Value = 13140.476
$ws.Range("L77").Value = 23283.125
$ws.Range("M77").Value = -8772.476000000001
$ws.Range("N77").Value = -32019.125

# Row 122 (G=36168)
$ws.Range("H122").Value = 3309.3845
$ws.Range("I122").Value = 1878
$ws.Range("J122").Value = 3945.5557
$ws.Range("K122").Value = 5634
$ws.Range("L122").Value = 11836.6671
$ws.Range("M122").Value = -3184
$ws.Range("N122").Value = -16736.6671

# Row 136 (G=43999)
$ws.Range("H136").Value = 1255560.1
$ws.Range("I136").Value = 3837.3333
$ws.Range("J136").Value = 2006593.8
$ws.Range("K136").Value = 11511.9999
$ws.Range("L136").Value = 6019781.4
$ws.Range("M136").Value = -8961.999899999999
$ws.Range("N136").Value = -6024881.4


$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G=44023)
$ws.Range("H31").Value = 4388.067
$ws.Range("I31").Value = 3881.3125
$ws.Range("J31").Value = 4967.2144
$ws.Range("K31").Value = 3881.3125
$ws.Range("L31").Value = 4967.2144
$ws.Range("M31").Value = -3586.3125
$ws.Range("N31").Value = -5557.2144

# Row 34 (G=44023)
$ws.Range("H34").Value = 4388.067
$ws.Range("I34").Value = 3881.3125
$ws.Range("J34").Value = 4967.2144
$ws.Range("K34").Value = 3881.3125
$ws.Range("L34").Value = 4967.2144
$ws.Range("M34").Value = -3679.3125
$ws.Range("N34").Value = -5371.2144

# Row 64 (G=10610)
$ws.Range("M64").ClearContents()
$ws.Range("H64").Value = 99998
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 99998
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 99998
$ws.Range("N64").Value = -100494

# Row 67 (G=10610)
$ws.Range("M67").ClearContents()
$ws.Range("H67").Value = 99998
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 99998
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 99998
$ws.Range("N67").Value = -101714

# Row 94 (G=32934)
$ws.Range("H94").Value = 1711.1765
$ws.Range("I94").Value = 1877
$ws.Range("J94").Value = 1563.7778
$ws.Range("K94").Value = 1877
$ws.Range("L94").Value = 1563.7778
$ws.Range("M94").Value = -1426
$ws.Range("N94").Value = -2465.7778

# Row 107 (G=27689)
$ws.Range("H107").Value = 1345.6666
$ws.Range("I107").Value = 1387.3077
$ws.Range("J107").Value = 1075
$ws.Range("K107").Value = 1387.3077
$ws.Range("L107").Value = 1075
$ws.Range("M107").Value = 532.6922999999999
$ws.Range("N107").Value = -4915

# Row 122 (G=36196)
$ws.Range("H122").Value = 1907
$ws.Range("I122").Value = 1284.4615
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 3853.3845
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -1403.3845
$ws.Range("N122").Value = -34900

# Row 132 (G=44019)
$ws.Range("H132").Value = 316200.44
$ws.Range("I132").Value = 3631.423
$ws.Range("J132").Value = 1670666.1
$ws.Range("K132").Value = 10894.269
$ws.Range("L132").Value = 5011998.300000001
$ws.Range("M132").Value = -8364.269
$ws.Range("N132").Value = -5017058.300000001


$ws = $wb.Worksheets.Item("CUL")
# Row 98 (G=19843)
$ws.Range("H98").Value = 667
$ws.Range("I98").Value = 559.75
$ws.Range("J98").Value = 728.2857
$ws.Range("K98").Value = 1679.25
$ws.Range("L98").Value = 2184.8571
$ws.Range("M98").Value = -181.25
$ws.Range("N98").Value = -5180.8571

# Row 112 (G=27855)
$ws.Range("H112").Value = 7399.9
$ws.Range("I112").Value = 2500
$ws.Range("J112").Value = 10666.5
$ws.Range("K112").Value = 7500
$ws.Range("L112").Value = 31999.5
$ws.Range("M112").Value = -6392
$ws.Range("N112").Value = -34215.5

# Row 113 (G=27843)
$ws.Range("M113").ClearContents()
$ws.Range("H113").Value = 3436.25
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3436.25
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 10308.75
$ws.Range("N113").Value = -14648.75

# Row 121 (G=27878)
$ws.Range("H121").Value = 9147.375
$ws.Range("I121").Value = 7475
$ws.Range("J121").Value = 11488.7
$ws.Range("K121").Value = 22425
$ws.Range("L121").Value = 34466.10000000001
$ws.Range("M121").Value = -21115
$ws.Range("N121").Value = -37086.10000000001


$ws = $wb.Worksheets.Item("GSM")
# Row 102 (G=36169)
$ws.Range("H102").Value = 2091.4897
$ws.Range("I102").Value = 851.55554
$ws.Range("J102").Value = 3613.2273
$ws.Range("K102").Value = 851.55554
$ws.Range("L102").Value = 3613.2273
$ws.Range("M102").Value = 770.44446
$ws.Range("N102").Value = -6857.2273

# Row 113 (G=27710)
$ws.Range("H113").Value = 3661.5454
$ws.Range("I113").Value = 1632.5
$ws.Range("J113").Value = 4112.4443
$ws.Range("K113").Value = 1632.5
$ws.Range("L113").Value = 4112.4443
$ws.Range("M113").Value = 537.5
$ws.Range("N113").Value = -8452.444299999999


$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("M22").ClearContents()
$ws.Range("H22").Value = 931.25
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 931.25
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 931.25
$ws.Range("N22").Value = -1521.25

# Row 27 (G=5277)
$ws.Range("M27").ClearContents()
$ws.Range("H27").Value = 931.25
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 931.25
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 931.25
$ws.Range("N27").Value = -1145.25

# Row 40 (G=36248)
$ws.Range("H40").Value = 4625.1875
$ws.Range("I40").Value = 3875.5
$ws.Range("J40").Value = 5374.875
$ws.Range("K40").Value = 3875.5
$ws.Range("L40").Value = 5374.875
$ws.Range("M40").Value = -3739.5
$ws.Range("N40").Value = -5646.875

# Row 43 (G=4314)
$ws.Range("H43").Value = 2504657
$ws.Range("I43").Value = 4169166.8
$ws.Range("J43").Value = 7892.25
$ws.Range("K43").Value = 4169166.8
$ws.Range("L43").Value = 7892.25
$ws.Range("M43").Value = -4168973.8
$ws.Range("N43").Value = -8278.25

# Row 46 (G=5282)
$ws.Range("H46").Value = 325991.72
$ws.Range("I46").Value = 3038
$ws.Range("J46").Value = 373836.72
$ws.Range("K46").Value = 3038
$ws.Range("L46").Value = 373836.72
$ws.Range("M46").Value = -2850
$ws.Range("N46").Value = -374212.72

# Row 93 (G=19993)
$ws.Range("H93").Value = 2543.9
$ws.Range("I93").Value = 2659.889
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 2659.889
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -1411.889
$ws.Range("N93").Value = -3996

# Row 101 (G=18549)
$ws.Range("H101").Value = 29544.8
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 29544.8
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 29544.8
$ws.Range("N101").Value = -36034.8

# Row 122 (G=36247)
$ws.Range("H122").Value = 9853
$ws.Range("I122").Value = 9907.385
$ws.Range("J122").Value = 9499.5
$ws.Range("K122").Value = 29722.155
$ws.Range("L122").Value = 28498.5
$ws.Range("M122").Value = -27272.155
$ws.Range("N122").Value = -33398.5


$ws = $wb.Worksheets.Item("WVR")
# Row 33 (G=2734)
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0

# Row 36 (G=2734)
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0

# Row 37 (G=3351)
$ws.Range("M37").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0

# Row 41 (G=21725)
$ws.Range("H41").Value = 25518
$ws.Range("I41").Value = 17500
$ws.Range("J41").Value = 26854.334
$ws.Range("K41").Value = 17500
$ws.Range("L41").Value = 26854.334
$ws.Range("M41").Value = -17110
$ws.Range("N41").Value = -27634.334

# Row 62 (G=12589)
$ws.Range("H62").Value = 2388227.5
$ws.Range("I62").Value = 5958069
$ws.Range("J62").Value = 8333.333000000001
$ws.Range("K62").Value = 5958069
$ws.Range("L62").Value = 8333.333000000001
$ws.Range("M62").Value = -5957445
$ws.Range("N62").Value = -9581.333000000001

# Row 64 (G=11036)
$ws.Range("H64").Value = 14973
$ws.Range("I64").Value = 14973
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 14973
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -14725

# Row 65 (G=12589)
$ws.Range("H65").Value = 2388227.5
$ws.Range("I65").Value = 5958069
$ws.Range("J65").Value = 8333.333000000001
$ws.Range("K65").Value = 29790345
$ws.Range("L65").Value = 41666.665
$ws.Range("M65").Value = -29787225
$ws.Range("N65").Value = -47906.665

# Row 67 (G=11036)
$ws.Range("H67").Value = 14973
$ws.Range("I67").Value = 14973
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 14973
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -14115

